$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-16 (columns A=aff_layers, B=optimizer, C=best_val_acc)
# to reflect the re-run grid search results.

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "nesterov"
$ws.Cells.Item(2, 3).Value = 0.3639968644810255

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "sgd"
$ws.Cells.Item(3, 3).Value = 0.3640276048631285

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "sgd"
$ws.Cells.Item(4, 3).Value = 0.448348472971519

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = "nesterov"
$ws.Cells.Item(5, 3).Value = 0.4528673091406526

$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "adamw"
$ws.Cells.Item(6, 3).Value = 0.3654262922488127

$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "sgd"
$ws.Cells.Item(7, 3).Value = 0.4578779914234334

$ws.Cells.Item(8, 1).Value = 2
$ws.Cells.Item(8, 2).Value = "adamw"
$ws.Cells.Item(8, 3).Value = 0.4527750879943438

$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "sgd"
$ws.Cells.Item(9, 3).Value = 0.4577089193218672

$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 2).Value = "nesterov"
$ws.Cells.Item(10, 3).Value = 0.4577857702771245

$ws.Cells.Item(11, 1).Value = 3
$ws.Cells.Item(11, 2).Value = "adamw"
$ws.Cells.Item(11, 3).Value = 0.4543428474815942

$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "sgd"
$ws.Cells.Item(12, 3).Value = 0.4535743379290205

$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 2).Value = "nesterov"
$ws.Cells.Item(13, 3).Value = 0.4568635588140361

$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = "nesterov"
$ws.Cells.Item(14, 3).Value = 0.4557107944851754

$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = "adamw"
$ws.Cells.Item(15, 3).Value = 0.4551882079894253

$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 2).Value = "adamw"
$ws.Cells.Item(16, 3).Value = 0.4569096693871905
